$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.848.74"
$ws.Range("E2").Value = "  -2.45%  "

$ws.Range("D3").Value = "2.417.71"
$ws.Range("E3").Value = "  -1.39%  "

$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.53%  "

$ws.Range("E7").Value = "  +0.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.526"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.33%  "

$ws.Range("D9").Value = "2.401.67"
$ws.Range("E9").Value = "  -1.94%  "

$ws.Range("E10").Value = "  -0.78%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.159"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.08"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.30%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.339"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.97%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.55%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000170"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.88%  "

$ws.Range("D16").Value = "2.844.67"
$ws.Range("E16").Value = "  -1.82%  "

$ws.Range("D17").Value = "60.733.46"
$ws.Range("E17").Value = "  -2.32%  "

$ws.Range("D18").Value = "2.404.68"
$ws.Range("E18").Value = "  -1.80%  "

$ws.Range("E19").Value = "  +5.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "322.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.72%  "

$ws.Range("E22").Value = "  -1.40%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.82%  "

$ws.Range("E24").Value = "  +0.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.55%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "64.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.26%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "582.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.94%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.46%  "

$ws.Range("D29").Value = "2.528.05"
$ws.Range("E29").Value = "  -1.80%  "

$ws.Range("D30").Value = "0.0₃0917"
$ws.Range("E30").Value = "  -4.41%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.44%  "

$ws.Range("E32").Value = "  -5.82%  "

$ws.Range("E33").Value = "  -2.09%  "

$ws.Range("E34").Value = "  -3.14%  "

$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("E36").Value = "  -5.64%  "

$ws.Range("E37").Value = "  -2.98%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "151.27"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.78%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.367"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.98%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.66%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.14"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.54%  "

$ws.Range("E42").Value = "  +0.06%  "

$ws.Range("E43").Value = "  -3.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.19"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.92%  "

$ws.Range("E46").Value = "  +13.91%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.80%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.61%  "

$ws.Range("E49").Value = "  -3.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0503"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.40%  "
